$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Tipo" column (currently column D)
$ws.Columns.Item(4).Insert()

# Header for the new column, matching the header style used by the other columns
$ws.Range("D1").Value = "MAE"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("D1").VerticalAlignment = -4160
$ws.Range("D1").Borders.LineStyle = 1

# Value for row 2
$ws.Range("D2").Value = 0.1401297014068613
